$wb = $excel.ActiveWorkbook

# --- Sheet 1: Site data ---
$ws1 = $wb.Worksheets.Item("Site data")
$ws1.Range("D4").Value = "Restore"
$ws1.Range("D5").Value = "Maintain"
$ws1.Range("D6").Value = "Maintain"
$ws1.Range("D7").Value = "Restore"

$ws1.Range("E4").Value = 172
$ws1.Range("F4").Value = 308
$ws1.Range("G4").Value = 575

$ws1.Range("E5").Value = 115
$ws1.Range("F5").Value = 255
$ws1.Range("G5").Value = 477

$ws1.Range("E6").Value = 196
$ws1.Range("F6").Value = 277
$ws1.Range("G6").Value = 562

$ws1.Range("E7").Value = 188
$ws1.Range("F7").Value = 316
$ws1.Range("G7").Value = 561

# --- Sheet 3: Feature data ---
$ws3 = $wb.Worksheets.Item("Feature data")
$ws3.Range("B4").Value = 6
$ws3.Range("C4").Value = 29

$ws3.Range("B5").Value = 76
$ws3.Range("C5").Value = 66

$ws3.Range("B6").Value = 84
$ws3.Range("C6").Value = 68

# --- Sheet 4: Expectation of "Maintain" ---
$ws4 = $wb.Worksheets.Item("Expectation of “Maintain”")
$ws4.Range("B4").Value = 6
$ws4.Range("C4").Value = 63
$ws4.Range("D4").Value = 21

$ws4.Range("B5").Value = 17
$ws4.Range("C5").Value = 64
$ws4.Range("D5").Value = 31

$ws4.Range("B6").Value = 83
$ws4.Range("C6").Value = 71
$ws4.Range("D6").Value = 35

$ws4.Range("B7").Value = 22
$ws4.Range("C7").Value = 68
$ws4.Range("D7").Value = 94

# --- Sheet 5: Expectation of "Restore" ---
$ws5 = $wb.Worksheets.Item("Expectation of “Restore”")
$ws5.Range("B4").Value = 34
$ws5.Range("C4").Value = 71
$ws5.Range("D4").Value = 6

$ws5.Range("B5").Value = 7
$ws5.Range("C5").Value = 20
$ws5.Range("D5").Value = 68

$ws5.Range("B6").Value = 28
$ws5.Range("C6").Value = 14
$ws5.Range("D6").Value = 29

$ws5.Range("B7").Value = 99
$ws5.Range("C7").Value = 30
$ws5.Range("D7").Value = 78

# --- Sheet 6: Expectation of "Signage" ---
$ws6 = $wb.Worksheets.Item("Expectation of “Signage”")
$ws6.Range("B4").Value = 70
$ws6.Range("C4").Value = 14
$ws6.Range("D4").Value = 83

$ws6.Range("B5").Value = 61
$ws6.Range("C5").Value = 91
$ws6.Range("D5").Value = 96

$ws6.Range("B6").Value = 82
$ws6.Range("C6").Value = 42
$ws6.Range("D6").Value = 19

$ws6.Range("B7").Value = 77
$ws6.Range("C7").Value = 12
$ws6.Range("D7").Value = 25
